$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column BC (56th... actually 55th) with nutrient "ethanol_g" / "ethanol" / unit "g"
$ws.Range("BC1").Value = "ethanol_g"
$ws.Range("BC2").Value = "ethanol"
$ws.Range("BC3").Value = "g"

# Match the style already used by the other row-2 header cells (e.g. D2/BB2)
$ws.Range("BC2").Style = $ws.Range("BB2").Style

# Update the view: scroll so column AH is the left-most visible column,
# and select BC3 as the active cell
$ws.Application.ActiveWindow.ScrollColumn = 34
$ws.Range("BC3").Select()
